$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.336.50"
$ws.Range("E2").Value = "  -3.15%  "
$ws.Range("D3").Value = "1.937.47"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.58"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7234"
$ws.Range("E6").Value = "  -6.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3339"
$ws.Range("E8").Value = "  -6.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.47"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07391"
$ws.Range("E10").Value = "  +4.70%  "
$ws.Range("E11").Value = "  -4.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08122"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.517"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "1.933.95"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.18"
$ws.Range("E15").Value = "  -6.18%  "
$ws.Range("E16").Value = "  -5.35%  "
$ws.Range("D17").Value = "30.328.16"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008354"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.64"
$ws.Range("E19").Value = "  -7.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.867"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "2.189.72"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.954"
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.809"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.49"
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.424"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1330"
$ws.Range("E29").Value = "  -9.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.560"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.445"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.246"
$ws.Range("E33").Value = "  -4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05201"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.255"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7511"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.736"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02001"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.842"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.673"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.66"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4547"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.027"
$ws.Range("E43").Value = "  -5.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8398"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.68"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.838"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.396"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.85"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.505"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4143"
$ws.Range("E51").Value = "  -4.91%  "
